$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45989, 11527.0059792127, 10332.3778751122, 16971.4, 6720.71718078781, 3.40396066249968),
    @(45990, 4520.63817582883, 7170.89237527922, 9003.4, 6143.57249361985, 179.627702870794),
    @(45991, 4158.11615452415, 6929.70487799966, 9003.4, 5880.39361914548, 158.612437381047),
    @(45992, 11309.5744899532, 10418.5344842323, 11232.26, 7172.62563212657, 264.954171514954),
    @(45993, 11230.0106043332, 10533.5655229883, 11232.26, 7110.20168839059, 267.146133807455),
    @(45994, 11775.3234136769, 10647.510122471, 11232.26, 7488.92967170283, 287.674158090577),
    @(45995, 11775.3234136769, 10455.1348357241, 11232.26, 7488.92967170283, 279.658521142788),
    @(45996, 11775.3234136769, 9722.35249155081, 11232.26, 7488.92967170283, 249.125923468902),
    @(45997, 4605.75571869631, 6433.63411109405, 11232.26, 7098.75248917023, 95.8386083443448),
    @(45998, 4450.60726673227, 6546.01555083205, 11232.26, 7083.68516255095, 99.8933630576249),
    @(45999, 4752.2157406543, 6887.10978072217, 11232.26, 7517.40805337532, 132.177409754062),
    @(46000, 12413.0231244078, 11022.7326710481, 11232.26, 7941.32071217855, 322.158057634445),
    @(46001, 12413.0231244078, 10991.4270666719, 11232.26, 7941.32071217855, 320.853657452103),
    @(46002, 12413.0231244078, 10974.8479592832, 11232.26, 7941.32071217855, 320.162861310907)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($row, $col).Value = $vals[$col - 1]
    }
}
